$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column B (Week_Start_Date) - shifts ASIN..is_holiday_week right by one
$ws.Columns.Item(2).Insert()

# Update header row
$ws.Range("B1").Value = "Week_Start_Date"

# Force column B to Text format so the date strings are stored as literal text, not date serials
$ws.Columns.Item(2).NumberFormat = "@"

# Update Week labels (remove leading zero) and set Week_Start_Date text values
$ws.Range("A2").Value = "W1"
$ws.Range("B2").Value = "2025-01-05"
$ws.Range("A3").Value = "W2"
$ws.Range("B3").Value = "2025-01-12"
$ws.Range("A4").Value = "W3"
$ws.Range("B4").Value = "2025-01-19"
$ws.Range("A5").Value = "W4"
$ws.Range("B5").Value = "2025-01-26"
$ws.Range("A6").Value = "W5"
$ws.Range("B6").Value = "2025-02-02"
$ws.Range("A7").Value = "W6"
$ws.Range("B7").Value = "2025-02-09"
$ws.Range("A8").Value = "W7"
$ws.Range("B8").Value = "2025-02-16"
$ws.Range("A9").Value = "W8"
$ws.Range("B9").Value = "2025-02-23"
$ws.Range("A10").Value = "W9"
$ws.Range("B10").Value = "2025-03-02"
$ws.Range("A11").Value = "W10"
$ws.Range("B11").Value = "2025-03-09"
$ws.Range("A12").Value = "W11"
$ws.Range("B12").Value = "2025-03-16"
$ws.Range("A13").Value = "W12"
$ws.Range("B13").Value = "2025-03-23"
$ws.Range("A14").Value = "W13"
$ws.Range("B14").Value = "2025-03-30"
$ws.Range("A15").Value = "W14"
$ws.Range("B15").Value = "2025-04-06"
$ws.Range("A16").Value = "W15"
$ws.Range("B16").Value = "2025-04-13"
$ws.Range("A17").Value = "W16"
$ws.Range("B17").Value = "2025-04-20"

# Update corrected MyForecast values (column D after insert)
$ws.Range("D5").Value = 63
$ws.Range("D8").Value = 75
$ws.Range("D9").Value = 71

# Convert is_holiday_week column (J) to boolean type
$ws.Range("J2:J17").Value = $false

# Update Summary sheet Max Forecast value (row 12) to reflect corrected forecast
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "75"

Write-Host "edit complete"
